$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.752.11'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.028.48'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.68'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.91'
$ws.Range('E6').Value = '  -3.13%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').Value = '3.028.91'
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.71'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('E14').Value = '  -3.73%  '
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '3.528.33'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.09'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '62.691.82'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '3.027.25'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '468.54'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.06'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.95'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E26').Value = '  -2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.37'
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.27'
$ws.Range('E29').Value = '  -2.30%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.47'
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.04'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '0.0₃0798'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.79'
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.27'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.04'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('E41').Value = '  -10.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '424.93'
$ws.Range('E42').Value = '  -2.35%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.114'
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.282'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').Value = '2.797.83'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '37.97'
$ws.Range('E47').Value = '  -7.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '128.80'
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.43'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('E51').Value = '  -0.24%  '
